# Applies corrected IFRS financial figures to rows 2-6 (years 2014-2018)
# and removes the stale/erroneous forecast rows 7-9 (2019E-2021E) data,
# keeping only their label columns (A, B, C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected values keyed by row number -> column letter -> value
$data = @{
    2 = @{ "D" = 7546; "E" = 267; "F" = 267; "G" = 11; "H" = 4; "I" = 22; "J" = -18; "K" = 11224; "L" = 8768; "M" = 2456; "N" = 1755; "O" = 701; "P" = 1217; "Q" = 89; "R" = -1444; "S" = 1519; "T" = 102; "U" = -13; "V" = 5326; "W" = 3.53; "X" = 0.05; "Y" = 1.17; "Z" = 0.05; "AA" = 356.99; "AB" = 88.64; "AC" = 91; "AD" = 263.95; "AE" = 7210; "AF" = 3.34; "AG" = 250; "AH" = 1.04; "AI" = 274.97; "AJ" = 24079554 }
    3 = @{ "D" = 12183; "E" = 985; "F" = 985; "G" = 567; "H" = 461; "I" = 310; "J" = 151; "K" = 12724; "L" = 9609; "M" = 3114; "N" = 2188; "O" = 926; "P" = 1217; "Q" = 1669; "R" = -863; "S" = 217; "T" = 626; "U" = 1044; "V" = 5935; "W" = 8.08; "X" = 3.78; "Y" = 15.73; "Z" = 3.85; "AA" = 308.55; "AB" = 123.46; "AC" = 1274; "AD" = 69.23; "AE" = 8989; "AF" = 9.81; "AG" = 250; "AH" = 0.28; "AI" = 19.67; "AJ" = 24079554 }
    4 = @{ "D" = 13008; "E" = 1269; "F" = 1269; "G" = 1074; "H" = 806; "I" = 578; "J" = 227; "K" = 13087; "L" = 9222; "M" = 3865; "N" = 2709; "O" = 1157; "P" = 1217; "Q" = 1100; "R" = -32; "S" = -1067; "T" = 414; "U" = 686; "V" = 5136; "W" = 9.75; "X" = 6.19; "Y" = 23.62; "Z" = 6.24; "AA" = 238.58; "AB" = 166.89; "AC" = 2376; "AD" = 25.13; "AE" = 11127; "AF" = 5.37; "AG" = 250; "AH" = 0.42; "AI" = 10.55; "AJ" = 24079554 }
    5 = @{ "D" = 9633; "E" = 880; "F" = 880; "G" = 888; "H" = 847; "I" = 684; "J" = 163; "K" = 11807; "L" = 6904; "M" = 4903; "N" = 3866; "O" = 1037; "P" = 1294; "Q" = 421; "R" = 78; "S" = -1550; "T" = 361; "U" = 60; "V" = 3782; "W" = 9.13; "X" = 8.79; "Y" = 20.8; "Z" = 6.8; "AA" = 140.82; "AB" = 187.62; "AC" = 2753; "AD" = 17.07; "AE" = 15049; "AF" = 3.12; "AG" = 500; "AH" = 1.06; "AI" = 18.8; "AJ" = 25618511 }
    6 = @{ "D" = 10263; "E" = 788; "F" = 788; "G" = 652; "H" = 636; "I" = 534; "K" = 11241; "L" = 5875; "M" = 5365; "N" = 4245; "P" = 1294; "Q" = 750; "R" = -364; "S" = -121; "T" = 865; "U" = -115; "V" = 3923; "W" = 7.68; "X" = 6.2; "Y" = 13.17; "Z" = 5.52; "AA" = 109.51; "AB" = 215.07; "AC" = 2064; "AD" = 13.25; "AE" = 16527; "AF" = 1.65; "AG" = 500; "AH" = 1.83; "AI" = 24.06; "AJ" = 25618511 }
}

foreach ($rowNum in $data.Keys) {
    $rowValues = $data[$rowNum]
    foreach ($col in $rowValues.Keys) {
        $ws.Range($col + $rowNum).Value = $rowValues[$col]
    }
}

# Rows 7-9 (2019E/2020E/2021E forecast columns) had all their figures
# removed in the source fix, leaving only the row index / period / ticker columns.
$ws.Range("D7:AJ9").ClearContents()
